# Refine data parser normalization, mapping, and logging
# Inserts additional "other item" / "balancing item" rows into the Income
# Statement template, expanding it from 55 to 70 data rows.
#
# Strategy: work from the BOTTOM of the sheet upward so that the row
# numbers referenced below always describe the *original* (pre-edit)
# layout and are never invalidated by an earlier insertion.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 8. Append a brand-new row after the old last row (55) for the new
#    "终止确认收益" line -> becomes row 70 in the final sheet.
$ws.Range("A56").Value = "以摊余成本计量的金融资产终止确认收益"

# 7. Insert 2 rows above old row 48 ("基本每股收益") -> final rows 60-61.
$ws.Range("A48:A49").EntireRow.Insert()
$ws.Range("A48").Value = "净利润其他项目"
$ws.Range("A49").Value = "净利润差额(合计平衡项目)"

# 6. Insert 1 row above old row 47 ("扣除非经常性损益后的净利润") -> final row 58.
$ws.Range("A47").EntireRow.Insert()
$ws.Range("A47").Value = "被合并方在合并前实现利润"

# 5. Insert 3 rows above old row 42 ("五、净利润") -> final rows 50-52.
$ws.Range("A42:A44").EntireRow.Insert()
$ws.Range("A42").Value = "未确认投资损失"
$ws.Range("A43").Value = "影响净利润的其他项目"
$ws.Range("A44").Value = "净利润差额(合计平衡项目2)"

# 4. Insert 2 rows above old row 40 ("四、利润总额") -> final rows 46-47.
$ws.Range("A40:A41").EntireRow.Insert()
$ws.Range("A40").Value = "影响利润总额的其他项目"
$ws.Range("A41").Value = "利润总额平衡项目"

# 3. Insert 4 rows above old row 35 ("三、营业利润") -> final rows 37-40.
$ws.Range("A35:A38").EntireRow.Insert()
$ws.Range("A35").Value = "资产减值损失(新)"
$ws.Range("A36").Value = "信用减值损失(新)"
$ws.Range("A37").Value = "营业利润其他项目"
$ws.Range("A38").Value = "营业利润平衡项目"

# 2. Insert 1 row above old row 28 ("加:其他收益") -> final row 29.
$ws.Range("A28").EntireRow.Insert()
$ws.Range("A28").Value = "营业总成本其他项目"

# 1. Insert 1 row above old row 9 ("二、营业总成本") -> final row 9.
$ws.Range("A9").EntireRow.Insert()
$ws.Range("A9").Value = "营业总收入其他项目"
